$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 16
$ws.Range("H2").Value = 16
$ws.Range("E4").Value = 11
$ws.Range("F5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("H7").Value = 2
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 2
$ws.Range("H8").Value = 2
$ws.Range("F15").Value = 51
$ws.Range("H15").Value = 51
$ws.Range("E17").Value = 69
$ws.Range("E18").Value = 62
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = 22
$ws.Range("E19").Value = 29
$ws.Range("F21").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("F24").Value = 7
$ws.Range("H24").Value = 7
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 4
$ws.Range("H25").Value = 4
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 4
$ws.Range("E33").Value = 19
$ws.Range("F35").Value = 2
$ws.Range("H35").Value = 2
$ws.Range("E36").Value = 55
$ws.Range("F37").Value = 12
$ws.Range("H37").Value = 12
$ws.Range("F39").Value = 8
$ws.Range("H39").Value = 8
$ws.Range("F40").Value = 5
$ws.Range("H40").Value = 5
$ws.Range("E42").Value = 22
$ws.Range("F42").Value = 6
$ws.Range("H42").Value = 6
$ws.Range("F45").Value = 9
$ws.Range("H45").Value = 9
$ws.Range("F47").Value = 22
$ws.Range("H47").Value = 22
$ws.Range("E49").Value = 41
$ws.Range("F49").Value = 21
$ws.Range("H49").Value = 21
$ws.Range("F52").Value = 1
$ws.Range("H52").Value = 1
$ws.Range("E59").Value = 7
$ws.Range("F60").Value = 5
$ws.Range("H60").Value = 5
$ws.Range("E61").Value = 18
$ws.Range("F61").Value = 4
$ws.Range("H61").Value = 4
$ws.Range("F62").Value = 3
$ws.Range("H62").Value = 3
$ws.Range("F63").Value = 3
$ws.Range("H63").Value = 3
$ws.Range("F66").Value = 9
$ws.Range("H66").Value = 9
$ws.Range("F67").Value = 15
$ws.Range("H67").Value = 15
$ws.Range("F69").Value = 5
$ws.Range("H69").Value = 5
$ws.Range("F72").Value = 12
$ws.Range("H72").Value = 12
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 5
$ws.Range("H73").Value = 5
$ws.Range("F76").Value = 9
$ws.Range("H76").Value = 9
$ws.Range("F77").Value = 8
$ws.Range("H77").Value = 8
$ws.Range("E78").Value = 16
$ws.Range("F79").Value = 5
$ws.Range("H79").Value = 5
$ws.Range("E89").Value = 18
$ws.Range("F89").Value = 7
$ws.Range("H89").Value = 7
